$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 91; $r++) {
    $cell = $ws.Range("E$r")
    if ($cell.Value2 -eq "FALSA") {
        $cell.Value2 = "FAKE"
    }
}
